$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clasificación")

# Update pair names for the first five rows with real player names
$ws.Range("A2").Value = "Jason/Jorge"
$ws.Range("A3").Value = "Alex/Keneth"
$ws.Range("A4").Value = "Teto/Pedro"
$ws.Range("A5").Value = "Keko/Memo"
$ws.Range("A6").Value = "Memin/Juan"

# Make "Clasificación" the active sheet with A7 selected
$ws.Activate()
$ws.Range("A7").Select()
